$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the time_taken (F) column timestamps on the "data" sheet ---
$ws.Range("F2").Value  = "2021-10-05 14:20:15.424236"
$ws.Range("F3").Value  = "2021-10-05 14:20:15.424244"
$ws.Range("F4").Value  = "2021-10-05 14:20:15.424247"
$ws.Range("F5").Value  = "2021-10-05 14:20:15.424250"
$ws.Range("F6").Value  = "2021-10-05 14:20:15.424253"
$ws.Range("F7").Value  = "2021-10-05 14:20:15.424255"
$ws.Range("F8").Value  = "2021-10-05 14:20:15.424258"
$ws.Range("F9").Value  = "2021-10-05 14:20:15.424260"
$ws.Range("F10").Value = "2021-10-05 14:20:15.424263"
$ws.Range("F11").Value = "2021-10-05 14:20:15.424266"
$ws.Range("F12").Value = "2021-10-05 14:20:15.424268"
$ws.Range("F13").Value = "2021-10-05 14:20:15.424271"
$ws.Range("F14").Value = "2021-10-05 14:20:15.424273"
$ws.Range("F15").Value = "2021-10-05 14:20:15.424276"
$ws.Range("F16").Value = "2021-10-05 14:20:15.424278"
$ws.Range("F17").Value = "2021-10-05 14:20:15.424281"
$ws.Range("F18").Value = "2021-10-05 14:20:15.424283"
$ws.Range("F19").Value = "2021-10-05 14:20:15.424286"
$ws.Range("F20").Value = "2021-10-05 14:20:15.424288"
$ws.Range("F21").Value = "2021-10-05 14:20:15.424291"
$ws.Range("F22").Value = "2021-10-05 14:20:15.424293"
$ws.Range("F23").Value = "2021-10-05 14:20:15.424296"
$ws.Range("F24").Value = "2021-10-05 14:20:15.424298"
$ws.Range("F25").Value = "2021-10-05 14:20:15.424301"
$ws.Range("F26").Value = "2021-10-05 14:20:15.424303"
$ws.Range("F27").Value = "2021-10-05 14:20:15.424306"
$ws.Range("F28").Value = "2021-10-05 14:20:15.424308"
$ws.Range("F29").Value = "2021-10-05 14:20:15.424311"
$ws.Range("F30").Value = "2021-10-05 14:20:15.424313"
$ws.Range("F31").Value = "2021-10-05 14:20:15.424316"
$ws.Range("F32").Value = "2021-10-05 14:20:15.424318"
$ws.Range("F33").Value = "2021-10-05 14:20:15.424320"
$ws.Range("F34").Value = "2021-10-05 14:20:15.424323"
$ws.Range("F35").Value = "2021-10-05 14:20:15.424326"
$ws.Range("F36").Value = "2021-10-05 14:20:15.424328"
$ws.Range("F37").Value = "2021-10-05 14:20:15.424331"
$ws.Range("F38").Value = "2021-10-05 14:20:15.424333"
$ws.Range("F39").Value = "2021-10-05 14:20:15.424336"
$ws.Range("F40").Value = "2021-10-05 14:20:15.424338"
$ws.Range("F41").Value = "2021-10-05 14:20:15.424340"
$ws.Range("F42").Value = "2021-10-05 14:20:15.424343"
$ws.Range("F43").Value = "2021-10-05 14:20:15.424346"
$ws.Range("F44").Value = "2021-10-05 14:20:15.424348"

# --- Add the new "metadata" worksheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Seed the header row (B1:F1) and A2 by copying from the "data" sheet so the
# bold/bordered/centered header style (and the numeric style used by A2) is
# reused instead of creating brand-new style entries. G1 reuses the same
# header style too (copied from F1, since "data" has no G column).
$ws.Range("B1:F1").Copy($meta.Range("B1:F1"))
$ws.Range("F1").Copy($meta.Range("G1"))
$ws.Range("A2").Copy($meta.Range("A2"))

# Header row values
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Familial hypercholesterolaemia"
$meta.Range("C2").Value = 6

# "1.28" must stay textual (not become the number 1.28). Build it as text in
# a scratch cell (using a leading apostrophe to force text entry) and paste
# just the value across, so the destination cell ends up as a plain,
# unstyled text cell (matching the source workbook).
$scratch = $ws.Range("Z1")
$scratch.Value = "'1.28"
$scratch.Copy()
$meta.Range("D2").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false

$meta.Range("E2").Value = "2021-07-01T12:36:20.933082Z"
$meta.Range("F2").Value = "2021-10-05 14:20:15.420551"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/6/?format=json"

# Keep the original sheet active/selected as before the edit
$ws.Activate()
